$wb = $excel.ActiveWorkbook

# --- Sheet1: fill in July-October 2022 statistics (rows 56-59, columns C:K) ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(35138, 14754, 1426, 118, 976, 2520, 150, 11990, 23148),
    @(40666, 16594, 1849, 158, 1075, 3082, 191, 14873, 25793),
    @(35124, 15371, 1649, 155, 997, 2801, 131, 14981, 20143),
    @(34049, 14400, 1437, 150, 1018, 2605, 210, 15020, 19029)
)

$startRow = 56
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $col = 3 + $c   # column C = 3
        $ws1.Cells.Item($row, $col).Value = $values[$c]
    }
}

# --- Sheet3: extend shared formula E3:E7 -> E3:E13, fill E8:E13 with =C-D ---
$ws3 = $wb.Worksheets.Item("Sheet3")
for ($r = 8; $r -le 13; $r++) {
    $ws3.Range("E" + $r).Formula = "=C" + $r + "-D" + $r
}

# --- Update sheet selections to match final state ---
$ws1.Activate()
$ws1.Range("C59:K59").Select()

$ws3.Activate()
$ws3.Range("E13").Select()
